$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00337C21" w:rsidRPr="00BB11F8" w:rsidRDefault="00BB11F8" w:rsidP="00337C21"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>jsp</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:include</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> page="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>header.jsp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00BB11F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>"/&gt;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="rIdHL1" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hipervnculo"/><w:lang w:val="en-US"/></w:rPr><w:t>https://wiki.jasig.org/display/CASC/Using+the+CAS+Client+3.1+with+Spring+Security</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="rIdHL2" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hipervnculo"/><w:lang w:val="en-US"/></w:rPr><w:t>http://java-assist.blogspot.com/2012/07/javalangnoclassdeffounderror_8541.html</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:hyperlink r:id="rIdHL3" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hipervnculo"/><w:lang w:val="en-US"/></w:rPr><w:t>http://cia.sourceforge.net/template-component-service/tattleTaleReport/jar/spring-expression-3.0.3.RELEASE.jar.html</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space="preserve">//posiblemente aumentar la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>version</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>http://stackoverflow.com/questions/13168215/propertytypedescriptor-and-spring-3-1-2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rIdHL1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://wiki.jasig.org/display/CASC/Using+the+CAS+Client+3.1+with+Spring+Security" TargetMode="External"/><Relationship Id="rIdHL2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://java-assist.blogspot.com/2012/07/javalangnoclassdeffounderror_8541.html" TargetMode="External"/><Relationship Id="rIdHL3" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://cia.sourceforge.net/template-component-service/tattleTaleReport/jar/spring-expression-3.0.3.RELEASE.jar.html" TargetMode="External"/></Relationships></pkg:xmlData></pkg:part></pkg:package>
"@
$r.InsertXML($xml) | Out-Null
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
